$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Slight wording change to the log key: the "direction" row now reads
# "direction (next/prev)" instead of "direction (forward/backward)".
$ws.Range("B10").Value = "direction (next/prev)"

# Reflect the last user selection being on B10 (as in the saved file).
$ws.Range("B10").Select()
